$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "27.482.15"
$ws.Cells.Item(2, 5).Value = "  +4.87%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.721.02"
$ws.Cells.Item(3, 5).Value = "  +4.11%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.10%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'229.07"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +4.51%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.5404"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +3.32%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'1.004"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.2752"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.19%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.06765"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +6.49%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'21.45"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +4.10%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.07799"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.63%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'4.703"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +2.82%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.735.68"
$ws.Cells.Item(13, 5).Value = "  -0.68%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.958.80"
$ws.Cells.Item(14, 5).Value = "  +3.99%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.5975"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +5.64%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.0₅8362"
$ws.Cells.Item(16, 5).Value = "  +2.03%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'68.57"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +4.70%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "27.460.36"
$ws.Cells.Item(18, 5).Value = "  +4.79%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'4.795"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.96%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.10%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'209.65"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +9.25%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'10.89"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +4.82%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'6.211"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +3.13%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'1.005"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.10%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'146.17"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.80%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.1248"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +3.70%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'7.428"
$ws.Cells.Item(27, 4).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 4).Value = "'16.80"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +5.16%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'1.621"
$ws.Cells.Item(29, 4).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 4).Value = "'0.05585"
$ws.Cells.Item(30, 4).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 4).Value = "'1.311"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.58%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.665"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +5.02%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'3.520"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +4.45%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.622"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.42%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.9742"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +3.08%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'2.853"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.83%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'2.443"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.47%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.5848"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.71%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.01645"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.91%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'5.839"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.13%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'1.003"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.08%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.040.80"
$ws.Cells.Item(42, 5).Value = "  +1.17%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'0.8372"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.15%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'102.44"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.43%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.864.26"
$ws.Cells.Item(45, 5).Value = "  +3.76%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "'59.49"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.77%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(47, 4).Value = "0.0₈111"
$ws.Cells.Item(47, 5).Value = "  +3.98%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'8.179"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.55%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.4434"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.87%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.9979"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.63%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.05274"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.73%  "
